# Migrazione dati sondaggio - aggiornamento 01.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1 (liste_naz): update percentages with new poll numbers ---
$ws1.Range("C2").Value = 0.238
$ws1.Range("F2").Value = 0.238

$ws1.Range("C3").Value = 0.221
$ws1.Range("F3").Value = 0.221

$ws1.Range("C4").Value = 0.14
$ws1.Range("F4").Value = 0.14

$ws1.Range("C5").Value = 0.112
$ws1.Range("F5").Value = 0.112

$ws1.Range("C6").Value = 0.074
$ws1.Range("F6").Value = 0.074

$ws1.Range("C7").Value = 0.049
$ws1.Range("F7").Value = 0.049

$ws1.Range("C8").Value = 0.038
$ws1.Range("F8").Value = 0.038

$ws1.Range("C9").Value = 0.027
$ws1.Range("F9").Value = 0.027

$ws1.Range("C10").Value = 0.025
$ws1.Range("F10").Value = 0.025

$ws1.Range("C11").Value = 0.023
$ws1.Range("F11").Value = 0.023

$ws1.Range("C12").Value = 0.01
$ws1.Range("F12").Value = 0.01
$ws1.Range("I12").Value = 202

# row 12 list name: "Altri di cdx" -> "Noi con l'Italia"
$ws1.Range("A12").Value = "Noi con l'Italia"

# row 13: becomes a generic "Altri 1" row, percentage updated, coalition col cleared
$ws1.Range("C13").Value = 0.01
$ws1.Range("F13").Value = 0.01
$ws1.Range("D13").ClearContents()
$ws1.Range("A13").Value = "Altri 1"

# new rows 14-17 (Altri 2..5) - set list names first so shared-string order
# matches "Altri 2/3/4/5" before the "ALTRI" coalition label is introduced
$ws1.Range("A14").Value = "Altri 2"
$ws1.Range("A15").Value = "Altri 3"
$ws1.Range("A16").Value = "Altri 4"
$ws1.Range("A17").Value = "Altri 5"

$ws1.Range("J13").Value = "ALTRI"

$ws1.Range("C14").Value = 0.01
$ws1.Range("C14").NumberFormat = $ws1.Range("C13").NumberFormat()
$ws1.Range("E14").Value = $false
$ws1.Range("F14").Value = 0.01
$ws1.Range("F14").NumberFormat = $ws1.Range("F13").NumberFormat()
$ws1.Range("G14").Value = 0.25
$ws1.Range("G14").NumberFormat = $ws1.Range("G13").NumberFormat()
$ws1.Range("H14").Value = $false
$ws1.Range("J14").Value = "ALTRI"

$ws1.Range("C15").Value = 0.01
$ws1.Range("C15").NumberFormat = $ws1.Range("C13").NumberFormat()
$ws1.Range("E15").Value = $false
$ws1.Range("F15").Value = 0.01
$ws1.Range("F15").NumberFormat = $ws1.Range("F13").NumberFormat()
$ws1.Range("G15").Value = 0.25
$ws1.Range("G15").NumberFormat = $ws1.Range("G13").NumberFormat()
$ws1.Range("H15").Value = $false
$ws1.Range("J15").Value = "ALTRI"

$ws1.Range("C16").Value = 0.01
$ws1.Range("C16").NumberFormat = $ws1.Range("C13").NumberFormat()
$ws1.Range("E16").Value = $false
$ws1.Range("F16").Value = 0.01
$ws1.Range("F16").NumberFormat = $ws1.Range("F13").NumberFormat()
$ws1.Range("G16").Value = 0.25
$ws1.Range("G16").NumberFormat = $ws1.Range("G13").NumberFormat()
$ws1.Range("H16").Value = $false
$ws1.Range("J16").Value = "ALTRI"

$ws1.Range("C17").Value = 0.003
$ws1.Range("C17").NumberFormat = $ws1.Range("C13").NumberFormat()
$ws1.Range("E17").Value = $false
$ws1.Range("F17").Value = 0.003
$ws1.Range("F17").NumberFormat = $ws1.Range("F13").NumberFormat()
$ws1.Range("G17").Value = 0.25
$ws1.Range("G17").NumberFormat = $ws1.Range("G13").NumberFormat()
$ws1.Range("H17").Value = $false
$ws1.Range("J17").Value = "ALTRI"

# --- sheet2 (altri_dati): add survey source column ---
$ws2.Columns.Item(1).Insert()
$ws2.Range("A1").Value = "Sondaggio"
$ws2.Range("A2").Value = "SWG 18/7/22"
$ws2.Range("B2").Value = 0.43
$ws2.Columns.Item(1).AutoFit()

# --- selections / active sheet ---
$ws1.Range("I14").Select()
$ws2.Activate()
$ws2.Range("B3").Select()
